# Applies "Fixed Stimulus Absolute Timestamps" edits:
#  - renames each worksheet tab (new timestamp suffixes)
#  - updates filename values in column B of each sheet (new timestamps,
#    in some cases new stim/category codes too)

$wb = $excel.ActiveWorkbook

# --- Rename worksheets (order matches sheetId 1..5) ---
$wsGNG  = $wb.Worksheets.Item(1)
$wsNB   = $wb.Worksheets.Item(2)
$wsRS   = $wb.Worksheets.Item(3)
$wsTOL  = $wb.Worksheets.Item(4)
$wsvSAT = $wb.Worksheets.Item(5)

$wsGNG.Name  = "GNG_TO-16504778469596746"
$wsNB.Name   = "NB_TO-16504778496806698"
$wsRS.Name   = "RS_TO-16504778496816685"
$wsTOL.Name  = "TOL_TO-16504778497456691"
$wsvSAT.Name = "vSAT_TO-16504778498086722"

# --- Sheet 1: GNG ---
$wsGNG.Range("B2").Value = "go_stims-16504778469216685.csv"
$wsGNG.Range("B3").Value = "GNG_stims-16504778469426727.csv"
$wsGNG.Range("B4").Value = "go_stims-16504778469446692.csv"
$wsGNG.Range("B5").Value = "GNG_stims-16504778469586725.csv"

# --- Sheet 2: NB ---
$wsNB.Range("B2").Value  = "OB-1650477847427706.csv"
$wsNB.Range("B3").Value  = "OB-1650477848038706.csv"
$wsNB.Range("B4").Value  = "ZB-match_2-16504778471686707.csv"
$wsNB.Range("B5").Value  = "OB-16504778477717052.csv"
$wsNB.Range("B6").Value  = "TB-1650477849663705.csv"
$wsNB.Range("B7").Value  = "ZB-match_1-16504778471396904.csv"
$wsNB.Range("B8").Value  = "TB-16504778493837047.csv"
$wsNB.Range("B9").Value  = "TB-16504778496246696.csv"
$wsNB.Range("B10").Value = "ZB-match_0-1650477847344673.csv"

# --- Sheet 3: RS ---
$wsRS.Range("B2").Value = "eyes closed"
$wsRS.Range("B3").Value = "eyes open"

# --- Sheet 4: TOL ---
$wsTOL.Range("B2").Value = "MM_stims-1650477849712706.csv"
$wsTOL.Range("B3").Value = "ZM_stims-1650477849688667.csv"
$wsTOL.Range("B4").Value = "MM_stims-16504778497286701.csv"
$wsTOL.Range("B5").Value = "ZM_stims-16504778497136676.csv"
$wsTOL.Range("B6").Value = "MM_stims-16504778497447064.csv"
$wsTOL.Range("B7").Value = "ZM_stims-16504778497296689.csv"

# --- Sheet 5: vSAT ---
$wsvSAT.Range("B2").Value = "vSAT_stims-1650477849792672.csv"
$wsvSAT.Range("B3").Value = "vSAT_stims-16504778497777128.csv"
$wsvSAT.Range("B4").Value = "SAT_stims-16504778497607026.csv"
$wsvSAT.Range("B5").Value = "SAT_stims-1650477849747671.csv"
